$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'333.28"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'1.34%"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'43.94"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'6.17%"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'5.786"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'3.14%"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'0.08338"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'2.10%"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'8.809"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'0.41%"
$c.Style = "Normal"
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$c = $ws.Range("D7")
$c.Value = "'4.503"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'-0.61%"
$c.Style = "Normal"
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Range("D8")
$c.Value = "'1.978"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'-2.58%"
$c.Style = "Normal"
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$c = $ws.Range("D9")
$c.Value = "'2.886"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'-1.92%"
$c.Style = "Normal"
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D10")
$c.Value = "'0.9340"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'1.69%"
$c.Style = "Normal"
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c = $ws.Range("D11")
$c.Value = "'0.1253"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'-1.40%"
$c.Style = "Normal"
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c = $ws.Range("D12")
$c.Value = "'0.1961"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'0.53%"
$c.Style = "Normal"
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c = $ws.Range("D13")
$c.Value = "'0.09507"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'1.79%"
$c.Style = "Normal"
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c = $ws.Range("D14")
$c.Value = "'0.03940"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'5.35%"
$c.Style = "Normal"
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c = $ws.Range("D15")
$c.Value = "'0.1067"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'0.54%"
$c.Style = "Normal"
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c = $ws.Range("D16")
$c.Value = "'0.001315"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'0.56%"
$c.Style = "Normal"
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$c = $ws.Range("D17")
$c.Value = "'0.006060"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'-2.57%"
$c.Style = "Normal"
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c = $ws.Range("D18")
$c.Value = "'3.506"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'1.94%"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'9.077"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'9.58%"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'-1.54%"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'0.2573"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'7.62%"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'0.04418"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'-0.02%"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'0.001260"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'-0.08%"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'1.51%"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'0.0001191"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'0.77%"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'0.0003996"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.02823"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'2.14%"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.05698"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'5.25%"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'0.007917"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'3.18%"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'0.89%"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'0.009004"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'0.06%"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'0.002157"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'-1.23%"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.01018"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'-10.90%"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'0.00007247"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'6.72%"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'-0.08%"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'0.003261"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'-6.83%"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'-0.12%"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'0.00002102"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'-0.08%"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'0.0002002"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'-0.08%"
$c.Style = "Normal"

Write-Host "Applied all changes"